$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.323.81"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").Value = "3.219.14"
$ws.Range("E3").Value = "  -1.55%  "
$ws.Range("E4").Value = "  +0.02%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "577.86"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -1.47%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "183.28"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -1.64%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.98%  "
$ws.Range("D9").Value = "3.217.08"
$ws.Range("E9").Value = "  -1.60%  "
$ws.Range("E10").Value = "  -3.04%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "6.56"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -2.30%  "
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("D13").Value = "3.778.30"
$ws.Range("E13").Value = "  -1.62%  "
$ws.Range("E14").Value = "  +0.11%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "27.61"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -3.35%  "
$ws.Range("D16").Value = "67.423.51"
$ws.Range("E16").Value = "  -0.81%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.0000168"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -2.25%  "
$ws.Range("D18").Value = "3.199.39"
$ws.Range("E18").Value = "  -2.34%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "5.73"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -2.15%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "13.41"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -1.74%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "395.10"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +3.41%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "7.52"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -2.53%  "
$ws.Range("E23").Value = "  +0.22%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "71.12"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -0.39%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "0.514"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("E26").Value = "  -3.13%  "
$ws.Range("E27").Value = "  -2.77%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "9.49"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -3.32%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("E30").Value = "  -2.55%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "5.53"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -5.70%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "22.53"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -1.45%  "
$ws.Range("E33").Value = "  -3.93%  "
$ws.Range("E34").Value = "  +0.00%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.24"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -2.72%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "159.94"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -1.42%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.46"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -4.91%  "
$ws.Range("E38").Value = "  +0.85%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "26.34"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -0.75%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.800"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -4.46%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "4.53"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -1.34%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "6.49"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -4.46%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "2.45"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -6.18%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.0681"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -1.82%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "40.56"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -1.94%  "
$ws.Range("D46").Value = "2.589.54"
$ws.Range("E46").Value = "  -2.26%  "
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "332.99"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -2.58%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "24.44"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -3.81%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.0276"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -2.79%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "6.25"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -0.39%  "
$ws.Range("E51").Value = "  -1.58%  "
